# Auto-generated edit script applying the cryptos.xlsx diff
# (price "Price" / "Volume(1h)" refresh + a 3-way and a 2-way row
# swap among coin rows 39-41 and 48-49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.958.08"
$ws.Range("E2").Value = "  +0.78%  "

$ws.Range("D3").Value = "2.598.61"
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'524.24"
$ws.Range("E5").Value = "  +3.39%  "

$ws.Range("D6").Value = "'154.70"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = "  +1.81%  "

$ws.Range("D9").Value = "'6.72"
$ws.Range("E9").Value = "  +2.56%  "

$ws.Range("E10").Value = "  +2.28%  "

$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("D13").Value = "3.053.22"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").Value = "60.970.43"

$ws.Range("D15").Value = "'21.69"
$ws.Range("E15").Value = "  +1.21%  "

$ws.Range("E16").Value = "  +1.02%  "

$ws.Range("D17").Value = "2.598.25"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("E18").Value = "  -0.99%  "

$ws.Range("D19").Value = "'353.56"
$ws.Range("E19").Value = "  +2.39%  "

$ws.Range("D20").Value = "'10.59"

$ws.Range("E21").Value = "  +1.75%  "

$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "'61.12"
$ws.Range("E23").Value = "  +1.94%  "

$ws.Range("E24").Value = "  +1.59%  "

$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").Value = "2.711.54"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").Value = "0.0₃0845"
$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "  +0.80%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").Value = "'6.34"
$ws.Range("E31").Value = "  +11.24%  "

$ws.Range("E32").Value = "  +0.18%  "

$ws.Range("E33").Value = "  +3.40%  "

$ws.Range("D34").Value = "'148.61"
$ws.Range("E34").Value = "  -3.19%  "

$ws.Range("D35").Value = "'4.19"
$ws.Range("E35").Value = "  +5.26%  "

$ws.Range("D36").Value = "'0.943"
$ws.Range("E36").Value = "  +9.44%  "

$ws.Range("E37").Value = "  +1.12%  "

$ws.Range("E38").Value = "  +2.32%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'0.850"
$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'3.80"
$ws.Range("E40").Value = "  +1.19%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'36.49"
$ws.Range("E41").Value = "  +1.63%  "

$ws.Range("D42").Value = "'287.70"
$ws.Range("E42").Value = "  -3.05%  "

$ws.Range("E43").Value = "  +1.58%  "

$ws.Range("E44").Value = "  +1.52%  "

$ws.Range("E45").Value = "  +0.65%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").Value = "'19.61"
$ws.Range("E47").Value = "  -0.97%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0238"
$ws.Range("E48").Value = "  +2.23%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.88"
$ws.Range("E49").Value = "  +0.66%  "

$ws.Range("D50").Value = "'10.32"
$ws.Range("E50").Value = "  +0.11%  "

$ws.Range("D51").Value = "'19.05"
$ws.Range("E51").Value = "  +8.23%  "
